# Fruta / hortaliza, semanal
# Insert two new weekly records (row 424 and 425) for
# Feria Lagunitas de Puerto Montt - Naranja - Valencia,
# shifting the existing rows 424-444 down to 426-446.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 424, pushing the
# existing data (old rows 424-444) down to new rows 426-446.
$ws.Range("A424:A425").EntireRow.Insert()

# New row 424: Valencia / Primera
$ws.Range("A424").Value = 4
$ws.Range("B424").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C424").Value = "Los Lagos"
$ws.Range("D424").Value = 44714
$ws.Range("E424").Value = 10
$ws.Range("F424").Value = "Fruta"
$ws.Range("G424").Value = 100102
$ws.Range("H424").Value = "Cítricos"
$ws.Range("I424").Value = 100102005
$ws.Range("J424").Value = "Naranja"
$ws.Range("K424").Value = "Valencia"
$ws.Range("L424").Value = "Primera"
$ws.Range("M424").Value = 600
$ws.Range("N424").Value = 16000
$ws.Range("O424").Value = 17000
$ws.Range("P424").Value = 16500
$ws.Range("Q424").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R424").Value = "Región de O'Higgins"
$ws.Range("S424").Value = 1100
$ws.Range("T424").Value = 15

# New row 425: Valencia / Segunda
$ws.Range("A425").Value = 4
$ws.Range("B425").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C425").Value = "Los Lagos"
$ws.Range("D425").Value = 44714
$ws.Range("E425").Value = 10
$ws.Range("F425").Value = "Fruta"
$ws.Range("G425").Value = 100102
$ws.Range("H425").Value = "Cítricos"
$ws.Range("I425").Value = 100102005
$ws.Range("J425").Value = "Naranja"
$ws.Range("K425").Value = "Valencia"
$ws.Range("L425").Value = "Segunda"
$ws.Range("M425").Value = 300
$ws.Range("N425").Value = 14000
$ws.Range("O425").Value = 14000
$ws.Range("P425").Value = 14000
$ws.Range("Q425").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R425").Value = "Región de O'Higgins"
$ws.Range("S425").Value = 933
$ws.Range("T425").Value = 15

# Make sure the date column keeps the same numeric (serial) style it had
# elsewhere in column D (rather than Insert's own format guess).
$ws.Range("D424").NumberFormat = $ws.Range("D423").NumberFormat
$ws.Range("D425").NumberFormat = $ws.Range("D423").NumberFormat
